# Generate Report for Handoff
#
# A new handoff round was generated for the still-pending ("low" priority)
# files. Their Priority is bumped to "ht" (matching the already-handed-off
# rows) and their "Latest Handoff Datetime" is refreshed to the new
# generation timestamp, for both the zh-cn and de-de localization sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = "2016-08-29 22:32:57"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($row in 4..7) {
    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = "2016-08-29 22:33:05"
}

# The "Overview" sheet's "Latest HO Xliff Generate Date" column shares its
# text with the de-de sheet's "Latest Handoff Datetime" for these rows
# (both previously read "2016-08-29 22:32:45"), so it moves in lockstep.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $overview.Range("G$row").Value = "2016-08-29 22:33:05"
}
